$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'56.986.11"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").Value = "'2.398.02"
$ws.Range("E3").Value = "  -3.64%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'486.10"
$ws.Range("E5").Value = "  -1.07%  "
$ws.Range("D6").Value = "'154.74"
$ws.Range("E6").Value = "  +1.83%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.598"
$ws.Range("E8").Value = "  +16.67%  "
$ws.Range("D9").Value = "'2.413.51"
$ws.Range("E9").Value = "  -3.37%  "
$ws.Range("D10").Value = "'6.33"
$ws.Range("E10").Value = "  +10.83%  "
$ws.Range("D11").Value = "'0.0993"
$ws.Range("E11").Value = "  +0.77%  "
$ws.Range("D12").Value = "'0.334"
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("E13").Value = "  +1.36%  "
$ws.Range("D14").Value = "'2.819.55"
$ws.Range("E14").Value = "  -3.50%  "
$ws.Range("D15").Value = "'56.947.46"
$ws.Range("E15").Value = "  +0.43%  "
$ws.Range("D16").Value = "'20.65"
$ws.Range("E16").Value = "  -2.71%  "
$ws.Range("E17").Value = "  -2.32%  "
$ws.Range("D18").Value = "'2.414.04"
$ws.Range("E18").Value = "  -3.38%  "
$ws.Range("D19").Value = "'4.72"
$ws.Range("E19").Value = "  +3.70%  "
$ws.Range("D20").Value = "'323.79"
$ws.Range("E20").Value = "  +0.71%  "
$ws.Range("D21").Value = "'9.89"
$ws.Range("E21").Value = "  -4.06%  "
$ws.Range("D22").Value = "'0.997"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").Value = "'5.95"
$ws.Range("E23").Value = "  +1.05%  "
$ws.Range("D24").Value = "'58.17"
$ws.Range("E24").Value = "  -0.34%  "
$ws.Range("D25").Value = "'0.404"
$ws.Range("E25").Value = "  -1.74%  "
$ws.Range("E26").Value = "  -0.34%  "
$ws.Range("D27").Value = "'0.160"
$ws.Range("E27").Value = "  -0.66%  "
$ws.Range("D28").Value = "'2.516.08"
$ws.Range("E28").Value = "  -2.97%  "
$ws.Range("D29").Value = "'7.23"
$ws.Range("E29").Value = "  -4.36%  "
$ws.Range("D30").Value = "'0.0₃0778"
$ws.Range("E30").Value = "  -2.99%  "
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("D32").Value = "'149.87"
$ws.Range("E32").Value = "  -1.03%  "
$ws.Range("D33").Value = "'18.49"
$ws.Range("E33").Value = "  +0.97%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").Value = "'5.26"
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("D36").Value = "'1.15"
$ws.Range("E36").Value = "  -0.50%  "
$ws.Range("E37").Value = "  -1.51%  "
$ws.Range("D38").Value = "'0.840"
$ws.Range("E38").Value = "  -3.42%  "
$ws.Range("D39").Value = "'34.04"
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("D40").Value = "'0.100"
$ws.Range("E40").Value = "  +8.04%  "
$ws.Range("E41").Value = "  -0.31%  "
$ws.Range("D42").Value = "'1.36"
$ws.Range("E42").Value = "  -1.81%  "
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Value = "'0.594"
$ws.Range("E44").Value = "  -3.20%  "
$ws.Range("D45").Value = "'268.02"
$ws.Range("E45").Value = "  +0.68%  "
$ws.Range("D46").Value = "'0.0528"
$ws.Range("E46").Value = "  -6.03%  "
$ws.Range("E47").Value = "  -0.33%  "
$ws.Range("E48").Value = "  -0.45%  "
$ws.Range("E49").Value = "  -5.60%  "
$ws.Range("D50").Value = "'1.870.72"
$ws.Range("E50").Value = "  -1.16%  "
$ws.Range("D51").Value = "'17.37"
$ws.Range("E51").Value = "  -2.40%  "
